$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ComponentsAndLicenses")

# Insert a new column before the existing "Normalized License Type" column (col H)
# for the new "PackageUrl" field; this shifts columns H:S to I:T.
$ws.Columns.Item(8).Insert()

# Match the width of the neighboring "Usage pattern" column (G) as closely as possible.
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth

# Populate the new column's header (row 1) and template placeholder (row 2).
# Row 2 is set first so the shared-string table order matches ($packageUrl$ before PackageUrl).
$ws.Range("H2").Value = "`$packageUrl`$"
$ws.Range("H1").Value = "PackageUrl"

# Update the active selection/view as left by the editor.
$ws.Range("G12").Select()
